$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "ClkViewScreen事件"

# Correct the "attribute display name" (E column) values so that they line up
# with their corresponding "attribute variable name" (D column) values.
$ws.Range("E7").Value = "页面地址"
$ws.Range("E8").Value = "页面路径"
$ws.Range("E9").Value = "页面标题"
$ws.Range("E10").Value = "是否首日访问"
$ws.Range("E11").Value = "是否首次触发事件"
$ws.Range("E12").Value = "向前域名"

# Add a new (empty, but formatted) row to the table for an additional
# statistical indicator description, matching the style used by the rest of
# the E4:E13 column.
$ws.Range("E13").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Extend the duplicate-value conditional formatting to cover the new row.
$cf = $ws.Range("E4:E13").FormatConditions.Item(1)
$cf.ModifyAppliesToRange($ws.Range("E4:E14"))

# Update the active selection to reflect where the author left off editing.
$ws.Range("D31").Select() | Out-Null
